$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8223676666666666
$ws.Range("H2").Value = 2.467103
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 45.977138744652
$ws.Range("R2").Value = 413.794248701868
$ws.Range("S2").Value = 0.6412441619121594
$ws.Range("T2").Value = 0.6412441619121594

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8223676666666666
$ws.Range("H3").Value = 2.467103
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 4.416670290542667
$ws.Range("R3").Value = 39.75003261488401
$ws.Range("S3").Value = 0.06159939735768789
$ws.Range("T3").Value = 0.06159939735768789

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8223676666666666
$ws.Range("H4").Value = 2.467103
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 21.30608544423478
$ws.Range("R4").Value = 191.754768998113
$ws.Range("S4").Value = 0.2971564407301527
$ws.Range("T4").Value = 0.2971564407301527
